# Added new sensors SHT30/31/35, GXHT30/31/35
# Inserts three new rows (12-14) above the existing LM75 row with spec
# data for the SHT30/GXHT30, SHT31/GXHT31 and SHT35/GXHT35 humidity &
# temperature sensors, shifting the rest of the sensor table down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before row 12 (LM75), pushing everything else down.
$ws.Rows("12:14").Insert()

# The freshly inserted rows come back with "no border" formatting, so copy
# the formatting (borders / number formats / fills) from the row just above
# (row 11, AHT20) down onto the new rows before filling in values.
$ws.Range("A11:L11").Copy()
$ws.Range("A12:L14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 12: SHT30/GXHT30 ---
$ws.Range("A12").Value = "SHT30/GXHT30"
$ws.Range("B12").Value = "2.2…5.5V"
$ws.Range("C12").Value = "I2C"
$ws.Range("D12").Value = "-40…125 °C"
$ws.Range("E12").Value = " ±0.3 °C"
$ws.Range("F12").Value = "0.06 °C"
$ws.Range("G12").Value = "0...100%"
$ws.Range("H12").Value = "±3%"
$ws.Range("I12").Value = "0.1 %"

# --- Row 13: SHT31/GXHT31 ---
$ws.Range("A13").Value = "SHT31/GXHT31"
$ws.Range("B13").Value = "2.2…5.5V"
$ws.Range("C13").Value = "I2C"
$ws.Range("D13").Value = "-40…125 °C"
$ws.Range("E13").Value = " ±0.2 °C"
$e13 = $ws.Range("E13").Characters(7, 2)
$e13.Font.Name = "Calibri"
$ws.Range("F13").Value = "0.06 °C"
$ws.Range("G13").Value = "0...100%"
$ws.Range("H13").Value = "±2%"
$ws.Range("I13").Value = "0.1 %"

# --- Row 14: SHT35/GXHT35 ---
$ws.Range("A14").Value = "SHT35/GXHT35"
$ws.Range("B14").Value = "2.2…5.5V"
$ws.Range("C14").Value = "I2C"
$ws.Range("D14").Value = "-40…125 °C"
$ws.Range("E14").Value = " ±0.1 °C"
$e14 = $ws.Range("E14").Characters(7, 2)
$e14.Font.Name = "Calibri"
$ws.Range("F14").Value = "0.06 °C"
$ws.Range("G14").Value = "0...100%"
$ws.Range("H14").Value = "±1.5%"
$ws.Range("I14").Value = "0.1 %"

# Restore the cursor/selection the way it ends up after this kind of edit.
$ws.Range("M21").Select()
